# [VM:timothy.queen@3/25/2015 1:47:08 PM] updated risks
# Applies the "updated risks" edit to the Integrated Register sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Hide rows 4 and 8 (risks no longer actively shown on the register) ---
$ws.Rows.Item(4).Hidden = $true
$ws.Rows.Item(8).Hidden = $true

# --- Row 11: bump Date Modified ---
$ws.Range("E11").Value = 42088

# --- Row 12: fill in the rest of the "access to warning information" risk ---
$ws.Range("E12").Value = 42088
$ws.Range("F12").Value = "Access to warnings must be limited to specific individuals"
$ws.Range("G12").Value = "John"
$ws.Range("H12").Value = "Team"
$ws.Range("I12").Value = "eCL team"
$ws.Range("J12").Value = "System"
$ws.Range("K12").Value = 0.2
$ws.Range("L12").Value = 3
$ws.Range("Q12").Value = "Access is controlled through an access control list, supervisors and managers"
$ws.Range("R12").Value = 42064
$ws.Range("T12").Value = "Ensure all staff on the project understand the sensitivity of the data"
$ws.Range("U12").Value = "Exposure"

# --- Row 11: update the Contingency Plan text (new shared string, added last) ---
$ws.Range("T11").Value = "individuals identified in OY3 will be more aligned with the project"

# Force the Risk Exposure Ranking formula on row 12 to recompute now that
# Probability/Impact are numeric instead of blank placeholders.
$f12 = $ws.Range("N12").Formula()
$ws.Range("N12").Formula = $f12

# --- Selection / view bookkeeping ---
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 3
$ws.Range("E6").Select() | Out-Null
